$wb = $excel.ActiveWorkbook

# "pubmedlink_existance" sheet: correct the expected count in E2 (2 -> 1)
$ws3 = $wb.Worksheets.Item("pubmedlink_existance")
$ws3.Range("E2").Value = 1
$ws3.Range("E2").Select()

# "pages_with_pdq_citations" sheet: correct the expected count in E2 (3 -> 2).
# Activate this sheet last so it stays the active/selected tab, matching
# where the author's cursor ended up after making the edits.
$ws1 = $wb.Worksheets.Item("pages_with_pdq_citations")
$ws1.Range("E2").Value = 2
$ws1.Activate()
$ws1.Range("E3").Select()
